# Author_form.xlsx - Sarah Wald
# "Added some random interests to correct profile for Sarah Wald"
#
# The "Interests" section (header in A35) has a few blank rows below it
# (A36, A37, ...) ready to be filled in. Add two interests.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A36").Value = "Hiking"
$ws.Range("A37").Value = "Outdoors"

# Clear the (invisible, "no fill") background formatting that had been
# explicitly applied to A27 - matches the author's session where that
# cell's formatting got reset to the sheet default while editing nearby.
$ws.Range("A27").Interior.Pattern = -4142

# Leave the selection where the author's cursor ended up after typing the
# two new interest rows (Enter moves the active cell down each time).
$ws.Range("A38").Select() | Out-Null
